$d = $word.ActiveDocument

# Update the date heading (unique text, safe to use Find/Replace).
$d.Content.Find.Execute("2025-09-28 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-29 Monday", 2)

# Update the division problems in the table. Several values repeat across
# the sheet (e.g. "38÷6=" and "16÷5=" each occur twice), so a blind
# document-wide Find/Replace would clobber the wrong cell. Instead, target
# each table cell directly by its row/column position and overwrite its
# text, which also preserves the existing run formatting (font/size).
$t = $d.Tables.Item(1)

function Set-Cell($rowIndex, $colIndex, $newText) {
    $cell = $t.Rows.Item($rowIndex).Cells.Item($colIndex)
    $cell.Range.Text = $newText
}

# Row 1
Set-Cell 1 1 "19÷6="
Set-Cell 1 2 "91÷2="
Set-Cell 1 3 "36÷3="
Set-Cell 1 4 "57÷6="
Set-Cell 1 5 "50÷2="

# Row 5
Set-Cell 5 1 "89÷3="
Set-Cell 5 2 "32÷3="
Set-Cell 5 3 "41÷6="
Set-Cell 5 4 "64÷6="
Set-Cell 5 5 "51÷7="

# Row 9
Set-Cell 9 1 "38÷7="
Set-Cell 9 2 "29÷3="
Set-Cell 9 3 "53÷6="
Set-Cell 9 4 "50÷2="
Set-Cell 9 5 "29÷4="

# Row 13
Set-Cell 13 1 "71÷5="
Set-Cell 13 2 "47÷6="
Set-Cell 13 3 "98÷4="
Set-Cell 13 4 "85÷2="
Set-Cell 13 5 "68÷9="

# Row 17
Set-Cell 17 1 "74÷6="
Set-Cell 17 2 "89÷2="
Set-Cell 17 3 "96÷4="
Set-Cell 17 4 "80÷4="
Set-Cell 17 5 "16÷5="
